$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.196438908576965
$ws.Range("B1").Value = 2.194975852966309
$ws.Range("C1").Value = 1.999509692192078
$ws.Range("D1").Value = 1.530771493911743
$ws.Range("E1").Value = 0.8582099676132202
